$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("E2").Value = 1251282873
$ws.Range("O2").Value = 617

# Update selection / view: select A2 and scroll so A1 is the top-left cell visible
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
